# Agreement testcases code added
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "RecordingInfo" worksheet as the last tab, by
#    copying the AgreementObligations sheet (same column layout/styles)
#    and replacing its content with the Recording Info test data.
# ------------------------------------------------------------------
$srcForCopy = $wb.Worksheets.Item("AgreementObligations")
$srcForCopy.Copy($null, $srcForCopy)
$recInfo = $wb.Worksheets.Item($wb.Worksheets.Count)
$recInfo.Name = "RecordingInfo"

$recInfo.Cells.ClearContents()

$recInfo.Range("A1:I1").Value = "String"

$recInfo.Range("A2").Value = "Title"
$recInfo.Range("B2").Value = "DocumentType"
$recInfo.Range("C2").Value = "Book"
$recInfo.Range("D2").Value = "Volume"
$recInfo.Range("E2").Value = "Page"
$recInfo.Range("F2").Value = "State"
$recInfo.Range("G2").Value = "Country"
$recInfo.Range("H2").Value = "Grantor"
$recInfo.Range("I2").Value = "Grantee"

# Extend the header styling (bold/yellow fill) from F2 to the new G:I columns
$recInfo.Range("F2").Copy()
$recInfo.Range("G2:I2").PasteSpecial(-4122)

$recInfo.Range("A3").Value = "AddRecordingInfoALT"
$recInfo.Range("A4").Value = "AddRecordingInfoRAW"
$recInfo.Range("A5").Value = "EditRecordingInfoALT"
$recInfo.Range("A6").Value = "EditRecordingInfoRAW"

$recInfo.Range("H1").Select()

# ------------------------------------------------------------------
# 2. PayeeInfo: update the sample parcel string used for the ALT test
#    case (row 3).
# ------------------------------------------------------------------
$payeeInfo = $wb.Worksheets.Item("PayeeInfo")
$payeeInfo.Range("C3").Value = "Parcel #: 06514, Grantor Name: , County PID: IND001"
$payeeInfo.Range("C4").Select()

# ------------------------------------------------------------------
# 3. AgreementForm: remove the DOT test-case row and repurpose the
#    remaining ROW row with the new sample value.
# ------------------------------------------------------------------
$agreementForm = $wb.Worksheets.Item("AgreementForm")
$agreementForm.Rows(4).Delete()
$agreementForm.Range("B4").Value = "Sample Lease Workflow"
$agreementForm.Range("B4").Select()

# Leave AgreementForm as the active sheet/tab, as in the source file.
$agreementForm.Select()
$wb.Windows.Item(1).ScrollWorkbookTabs(2)
